$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "project" / "Description" (A1 previously held the "Juniper Nursing Home"
# index that is being replaced below, so re-point it at "project")
$ws.Range("A1").Value = "project"

# Row 2: rename the single "Juniper Nursing Home" entry to the first of three
$ws.Range("B2").Value = "Juniper Nursing Home 1"

# Row 3: new project row
$ws.Range("A3").Value = 1345007
$ws.Range("B3").Value = "Juniper Nursing Home 2"
$ws.Range("A3").HorizontalAlignment = -4131

# Row 4: new project row
$ws.Range("A4").Value = 1345008
$ws.Range("B4").Value = "Juniper Nursing Home 3"
$ws.Range("A4").HorizontalAlignment = -4131

# Move the active selection below the new data, as in the saved workbook
$ws.Range("A5").Select() | Out-Null
